$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous (long) query text, now kept verbatim but re-quoted/escaped
# as literal text (as if pasted back in from a quoted CSV export).
$oldQueryQuoted = """SELECT DISTINCT`n    ds.dataset_title AS """"Title"""",`n    ds.dataset_source_id AS """"Source ID"""", `n    ds.primary_disease AS """"Primary Disease"""",`n    CAST(ds.participant_count AS INT) AS """"Participants Count"""",`n    CAST(ds.sample_count AS INT) AS """"Sample Count"""",`n    CASE `n        WHEN LENGTH(TRIM(ds.description)) > 500 THEN`n            CASE`n                WHEN SUBSTR(TRIM(ds.description), 500, 1) = ' '`n                    THEN SUBSTR(TRIM(ds.description), 1, 499) || ' ...'`n                ELSE SUBSTR(TRIM(ds.description), 1, 500) || ' ...'`n            END`n        ELSE TRIM(ds.description)`n    END AS """"Description""""`nFROM df_dbgap ds`nORDER BY CAST(ds.dataset_title AS TEXT) ASC;"""

# New, simplified query that replaces it as the "live" DatasetsTab query.
$newQuery = "SELECT DISTINCT`n    REPLACE(ds.dataset_title, '  ', ' ') AS ""Title"",`n    ds.dataset_source_id AS ""Source ID"", `n    ds.primary_disease AS ""Primary Disease"",`n    -- CAST(ds.participant_count AS INT) AS ""Participants Count"",`n    CAST(ds.sample_count AS INT) AS ""Sample Count""`nFROM df_dbgap ds`nORDER BY CAST(ds.dataset_title AS TEXT) ASC;"

# Write the quoted/escaped text into the old query's cell first, then copy
# it down to B3, and only then overwrite B2 with the new simplified query.
$ws.Range("B2").Value = $oldQueryQuoted
$ws.Range("B3").Value = $oldQueryQuoted
$ws.Range("B2").Value = $newQuery

# Selection moves from B2 to C2 (and the frozen top-left scroll position
# on the sheet view is cleared).
$null = $ws.Range("C2").Select()
